# "Penalty Reward System" attempt (unfinished, per commit message):
# the author removed one row of data from each sheet by deleting the
# row in-place (so subsequent rows shift up), rather than clearing values.
#
# Sheet "Weekly Quantity": rows 5 and 6 (45361.99999999999/48 and
# 45368.99999999999/42) are removed; the old rows 7-8 shift up to become
# the new rows 5-6. Net effect: delete rows 5:6.
#
# Sheet "Monthly Trend": row 4 (45382.99999999999/90) is removed; the old
# row 5 shifts up to become the new row 4. Net effect: delete row 4.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows("5:6").Delete()

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Rows("4:4").Delete()
